$d = $word.ActiveDocument

# --- Edit 1 ---------------------------------------------------------------
# "The WIPO and Gallup indexes are reported for about 140-150 countries"
# becomes "... about 140 countries" (drop the "-150" range).
$d.Content.Find.Execute(
    "140-150 countries",
    $false, $true, $false, $false, $false,
    $true, 1, $false,
    "140 countries",
    2
) | Out-Null

# --- Edit 2 ---------------------------------------------------------------
# Rewrite of the PCA paragraph:
#   "... innovativeness of countries in a statistically ..."
#   -> "... innovativeness of countries and variance in a statistically ..."
$d.Content.Find.Execute(
    "innovativeness of countries in a statistically significant manner.",
    $false, $true, $false, $false, $false,
    $true, 1, $false,
    "innovativeness of countries and variance in a statistically significant manner.",
    2
) | Out-Null

#   "... As I will only use PCA (with one dependent variable, innovativeness), I will have to
#    extract the relative weight of happiness from its contribution to the eigenvectors.
#    It may end up being a qualitative estimate. "
#   -> "... As I will only use PCA, I will have to extract the relative weight of happiness and
#    other variables from their contribution to the eigenvectors. "
$d.Content.Find.Execute(
    "use PCA (with one dependent variable, innovativeness), I will have to extract the relative weight of happiness from its contribution to the eigenvectors. It may end up being a qualitative estimate. ",
    $false, $true, $false, $false, $false,
    $true, 1, $false,
    "use PCA, I will have to extract the relative weight of happiness and other variables from their contribution to the eigenvectors. ",
    2
) | Out-Null
